# Every monthly worksheet (January..December) gets the same edit:
#   Row 2:  "Gas" / 90 / 150   ->  "Botox" / 0 / 0
#   Row 3:  "Rent" / 0 / 0     ->  "Gas" / 0 / 0      (values unchanged, only the label)
#   Row 4:  "Food" / 0 / 0     ->  unchanged, left untouched
#   Row 5:  (did not exist)    ->  "God" / 0 / 0       (brand-new row)
#
# Row 3's B/C cells already hold the text "0" - we reuse them (via Copy /
# PasteSpecial -Values) as the source for the other "0" cells so the new
# cells end up with the same cell formatting/type as their neighbours
# instead of whatever a freshly typed literal would produce.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- Row 2: Gas/90/150 -> Botox/0/0 ---
    $ws.Range("A2").Value = "Botox"

    $ws.Range("B3").Copy()
    $ws.Range("B2").PasteSpecial(-4163)   # xlPasteValues

    $ws.Range("C3").Copy()
    $ws.Range("C2").PasteSpecial(-4163)   # xlPasteValues

    # --- Row 3: Rent -> Gas (values B3/C3 stay 0/0, untouched) ---
    $ws.Range("A3").Value = "Gas"

    # --- Row 4: Food/0/0 left completely untouched ---

    # --- Row 5: new row "God"/0/0 ---
    # Copy formatting from row 4 so the new cells match the existing
    # label/value styles (s=1 for labels, s=2 for values).
    $ws.Range("A4").Copy()
    $ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("B4:C4").Copy()
    $ws.Range("B5:C5").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("A5").Value = "God"

    $ws.Range("B4").Copy()
    $ws.Range("B5").PasteSpecial(-4163)   # xlPasteValues

    $ws.Range("C4").Copy()
    $ws.Range("C5").PasteSpecial(-4163)   # xlPasteValues

    $excel.CutCopyMode = 0
}
